# WS_holdings.xlsx update:
#  - bump the "as of" date in the confidential disclaimer text (A16)
#  - refresh the Weight / Percent Change figures in D2:E13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so the cells can be written, then
# restore protection once the edits are done.
$ws.Unprotect()

$ws.Range("D2").Value = 0.03119933210182017
$ws.Range("E2").Value = 0.004048582995951566

$ws.Range("D3").Value = 0.02343557163353451
$ws.Range("E3").Value = 0.0007115749525616888

$ws.Range("D4").Value = 0.05233193706911774
$ws.Range("E4").Value = -0.001613646841862715

$ws.Range("D5").Value = 0.138187225087198
$ws.Range("E5").Value = -0.003720478809446703

$ws.Range("D6").Value = 0.03091903059201217
$ws.Range("E6").Value = 0.01144492131616603

$ws.Range("D7").Value = 0.1166182013134905
$ws.Range("E7").Value = 0.01300875837197335

$ws.Range("D8").Value = 0.1016373274563844
$ws.Range("E8").Value = 0.004228718514432517

$ws.Range("D9").Value = 0.02917583905063652
$ws.Range("E9").Value = 0.002480876576390401

$ws.Range("D10").Value = 0.1258400027154948
$ws.Range("E10").Value = 0.005075187969924899

$ws.Range("D11").Value = 0.2470958339034965
$ws.Range("E11").Value = 0.001364132411786345

$ws.Range("D12").Value = 0.1035596990768146
$ws.Range("E12").Value = 0.006139677666922472

$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 0.003529078797958274

# Disclaimer paragraph: "as of 2021-05-25" -> "as of 2021-05-26"
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-26 for illustrative purposes only and are subject to change."
$ws.Range("A16").Value = $disclaimer

# Setting a multi-line value marks the row with an explicit custom height;
# AutoFit puts it back to the sheet's default (matching the original file).
$ws.Rows(16).AutoFit()

$ws.Protect()
